$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 2 + 2 =
$ws.Range("A1").Value = "2 + 2 ="
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 6
$ws.Range("D1").Value = 1
$ws.Range("E1").Value = -4
$ws.Range("F1").Value = 4

# Row 2: 2 * 2 =
$ws.Range("A2").Value = "2 * 2 ="
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = -8
$ws.Range("F2").Value = 4

# Row 3: 2 - 2 =
$ws.Range("A3").Value = "2 - 2 ="
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 14
$ws.Range("F3").Value = 0

# Row 4: 4 - 3 =
$ws.Range("A4").Value = "4 - 3 = "
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = -1
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1

$ws.Range("A3:A4").NumberFormat = "d-mmm"

$ws.Range("H6").Select()
